# CS133JS Lab04 Rubric - "Updated instructions and new buggy betas"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix a typo in the Part 2 rubric line (dropped the stray "oop") ---
$ws.Range("A13").Value = "  Loop +  i/o  in the HTML file"

# --- Remove the underline styling that had been applied down the Possible/
#     Score columns for the individual criteria rows (D6:E18) ---
$ws.Range("D6:E18").Font.Underline = $false

# --- Append the new "buggy betas" rubric sections below the existing table ---
# (typed in this order: the three bug-description lines for each lettered
# part, then the three section headers)
$ws.Range("A23").Value = "Should have for loop"
$ws.Range("A27").Value = "Validation with while loop"
$ws.Range("A31").Value = "Repeat with do…while, validation with while"
$ws.Range("A28").Value = "No if for alert"
$ws.Range("A32").Value = "No if"
$ws.Range("B28").Value = "Check num range, y/n"
$ws.Range("A22").Value = "1) A. Countdown, B. Loan Repayment, C. Cycling Plan"
$ws.Range("A26").Value = "2) A. Kindergarten Admit, B. Beverage Labeling,  C. Trail Difficulty Rating"
$ws.Range("A30").Value = "3) A. Grade Level, B. State Tax, C. Trail Types"

# --- Scroll the view down to the newly-added content and leave the cursor
#     where editing left off ---
$win = $excel.ActiveWindow
$win.ScrollRow = 17
$win.ScrollColumn = 1
$ws.Range("G13").Select() | Out-Null
